$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 426, shifting the existing
# data (old rows 426-448) down to rows 428-450.
$ws.Rows("426:427").Insert()

# Populate the two newly inserted rows with their data.
# Row 426 (new)
$ws.Cells.Item(426, 1).Value = 10
$ws.Cells.Item(426, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(426, 3).Value = "La Araucanía"
$ws.Cells.Item(426, 4).Value = [DateTime]"2023-06-29"
$ws.Cells.Item(426, 5).Value = 9
$ws.Cells.Item(426, 6).Value = "Fruta"
$ws.Cells.Item(426, 7).Value = 100102
$ws.Cells.Item(426, 8).Value = "Cítricos"
$ws.Cells.Item(426, 9).Value = 100102006
$ws.Cells.Item(426, 10).Value = "Pomelo"
$ws.Cells.Item(426, 11).Value = "Start Ruby"
$ws.Cells.Item(426, 12).Value = "Primera"
$ws.Cells.Item(426, 13).Value = 140
$ws.Cells.Item(426, 14).Value = 15000
$ws.Cells.Item(426, 15).Value = 15000
$ws.Cells.Item(426, 16).Value = 15000
$ws.Cells.Item(426, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(426, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(426, 19).Value = 1000
$ws.Cells.Item(426, 20).Value = 15

# Row 427 (new)
$ws.Cells.Item(427, 1).Value = 10
$ws.Cells.Item(427, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(427, 3).Value = "La Araucanía"
$ws.Cells.Item(427, 4).Value = [DateTime]"2023-06-29"
$ws.Cells.Item(427, 5).Value = 9
$ws.Cells.Item(427, 6).Value = "Fruta"
$ws.Cells.Item(427, 7).Value = 100102
$ws.Cells.Item(427, 8).Value = "Cítricos"
$ws.Cells.Item(427, 9).Value = 100102006
$ws.Cells.Item(427, 10).Value = "Pomelo"
$ws.Cells.Item(427, 11).Value = "Start Ruby"
$ws.Cells.Item(427, 12).Value = "Primera"
$ws.Cells.Item(427, 13).Value = 100
$ws.Cells.Item(427, 14).Value = 15000
$ws.Cells.Item(427, 15).Value = 15000
$ws.Cells.Item(427, 16).Value = 15000
$ws.Cells.Item(427, 17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(427, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(427, 19).Value = 1071
$ws.Cells.Item(427, 20).Value = 14
